# Remove the empty Subtitle placeholder ("Subtitle 2") from the first
# slide (the title slide). The title shape ("Title 1" / "Loan Risk
# Analysis") is left untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($shp in $s.Shapes) {
    if ($shp.Name -eq "Subtitle 2") {
        $shp.Delete()
    }
}
